# Add new company rows (Sony, Nike, The Hartford, Keystone, EnsoraHealth)
# to the "Company list" tracker, matching the existing Company name / Careers
# site layout, and tidy up the stray cell formatting that had accumulated
# on the blank trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New company rows -------------------------------------------------
# Column B (careers site URL) is written before column A (company name)
# for each row, matching how the rest of the sheet's shared strings were
# originally authored.

$ws.Range("B29").Value = "https://sonyglobal.wd1.myworkdayjobs.com/en-US/SonyGlobalCareers/"
$ws.Range("A29").Value = "Sony"

$ws.Range("B30").Value = "https://nike.wd1.myworkdayjobs.com/nke/"
$ws.Range("A30").Value = "Nike"

$ws.Range("B31").Value = "https://thehartford.wd5.myworkdayjobs.com/en-US/Careers_External"
$ws.Range("A31").Value = "The Hartford"

$ws.Range("B32").Value = "https://keystone.wd5.myworkdayjobs.com/Keystone/"
$ws.Range("A32").Value = "Keystone"

$ws.Range("B33").Value = "https://therapybrands.wd1.myworkdayjobs.com/EnsoraHealth"
$ws.Range("A33").Value = "EnsoraHealth"

# --- Formatting cleanup -------------------------------------------------
# A15:A28 had a redundant "applyFont" style lingering on them; clear it so
# they fall back to the sheet's default (Normal) style.
$ws.Range("A15:A32").ClearFormats()

# --- Selection / scroll position ----------------------------------------
$ws.Activate()
$ws.Range("B36").Select()
